$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (participant id): change all 12 -> 14 for rows 2 through 33
# (values are stored as text, not numbers, so force a text number format)
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "14"
}

# Column H (path): insert "VR/" into instructions_videos paths for rows
# whose path starts with "./instructions_videos/"
$instructionRows = 2,4,6,8,10,11,12,14,16,17,18,20,22,23,24,26,28,29,30,32,33
foreach ($r in $instructionRows) {
    $cell = $ws.Cells.Item($r, 8)
    $current = $cell.Value()
    $cell.Value = $current -replace '^\./instructions_videos/', './instructions_videos/VR/'
}

# Swap video_id / exp video path pairs (video_id column stored as text):
# row 7: 9 -> 8 ; row 13: 8 -> 9
$e7 = $ws.Cells.Item(7, 5)
$e7.NumberFormat = "@"
$e7.Value = "8"
$ws.Cells.Item(7, 8).Value = "../stimuli/exp_videos/VR/8.mp4"

$e13 = $ws.Cells.Item(13, 5)
$e13.NumberFormat = "@"
$e13.Value = "9"
$ws.Cells.Item(13, 8).Value = "../stimuli/exp_videos/VR/9.mp4"

# row 19: 7 -> 4 ; row 25: 4 -> 7
$e19 = $ws.Cells.Item(19, 5)
$e19.NumberFormat = "@"
$e19.Value = "4"
$ws.Cells.Item(19, 8).Value = "../stimuli/exp_videos/VR/4.mp4"

$e25 = $ws.Cells.Item(25, 5)
$e25.NumberFormat = "@"
$e25.Value = "7"
$ws.Cells.Item(25, 8).Value = "../stimuli/exp_videos/VR/7.mp4"
